$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Hspg2"
$ws.Cells.Item(2, 3).Value = "Col13a1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 139.028825
$ws.Cells.Item(2, 8).Value = 417.086475
$ws.Cells.Item(2, 9).Value = 0.428090758569687
$ws.Cells.Item(2, 10).Value = 0.428090758569687
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3943663333333334
$ws.Cells.Item(2, 14).Value = 1.183099
$ws.Cells.Item(2, 15).Value = 0.393328993667368
$ws.Cells.Item(2, 16).Value = 0.393328993667368
$ws.Cells.Item(2, 17).Value = 54.82828794289168
$ws.Cells.Item(2, 18).Value = 493.454591486025
$ws.Cells.Item(2, 19).Value = 0.1683805072665152
$ws.Cells.Item(2, 20).Value = 0.1683805072665152
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Hspg2"
$ws.Cells.Item(3, 3).Value = "Col13a1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 139.028825
$ws.Cells.Item(3, 8).Value = 417.086475
$ws.Cells.Item(3, 9).Value = 0.428090758569687
$ws.Cells.Item(3, 10).Value = 0.428090758569687
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.043492
$ws.Cells.Item(3, 14).Value = 0.130476
$ws.Cells.Item(3, 15).Value = 0.04337759881273123
$ws.Cells.Item(3, 16).Value = 0.04337759881273123
$ws.Cells.Item(3, 17).Value = 6.046641656900001
$ws.Cells.Item(3, 18).Value = 54.4197749121
$ws.Cells.Item(3, 19).Value = 0.01856954918067366
$ws.Cells.Item(3, 20).Value = 0.01856954918067366
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Hspg2"
$ws.Cells.Item(4, 3).Value = "Col13a1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 139.028825
$ws.Cells.Item(4, 8).Value = 417.086475
$ws.Cells.Item(4, 9).Value = 0.428090758569687
$ws.Cells.Item(4, 10).Value = 0.428090758569687
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.564779
$ws.Cells.Item(4, 14).Value = 1.694337
$ws.Cells.Item(4, 15).Value = 0.5632934075199009
$ws.Cells.Item(4, 16).Value = 0.5632934075199009
$ws.Cells.Item(4, 17).Value = 78.52056075467502
$ws.Cells.Item(4, 18).Value = 706.685046792075
$ws.Cells.Item(4, 19).Value = 0.2411407021224982
$ws.Cells.Item(4, 20).Value = 0.2411407021224982
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Hspg2"
$ws.Cells.Item(5, 3).Value = "Col13a1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 105.0686953333333
$ws.Cells.Item(5, 8).Value = 315.206086
$ws.Cells.Item(5, 9).Value = 0.3235223881606854
$ws.Cells.Item(5, 10).Value = 0.3235223881606855
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3943663333333334
$ws.Cells.Item(5, 14).Value = 1.183099
$ws.Cells.Item(5, 15).Value = 0.393328993667368
$ws.Cells.Item(5, 16).Value = 0.393328993667368
$ws.Cells.Item(5, 17).Value = 41.43555612672378
$ws.Cells.Item(5, 18).Value = 372.9200051405141
$ws.Cells.Item(5, 19).Value = 0.127250735364106
$ws.Cells.Item(5, 20).Value = 0.127250735364106
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Hspg2"
$ws.Cells.Item(6, 3).Value = "Col13a1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 105.0686953333333
$ws.Cells.Item(6, 8).Value = 315.206086
$ws.Cells.Item(6, 9).Value = 0.3235223881606854
$ws.Cells.Item(6, 10).Value = 0.3235223881606855
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.043492
$ws.Cells.Item(6, 14).Value = 0.130476
$ws.Cells.Item(6, 15).Value = 0.04337759881273123
$ws.Cells.Item(6, 16).Value = 0.04337759881273123
$ws.Cells.Item(6, 17).Value = 4.569647697437333
$ws.Cells.Item(6, 18).Value = 41.12682927693601
$ws.Cells.Item(6, 19).Value = 0.01403362436057092
$ws.Cells.Item(6, 20).Value = 0.01403362436057092
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Hspg2"
$ws.Cells.Item(7, 3).Value = "Col13a1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 105.0686953333333
$ws.Cells.Item(7, 8).Value = 315.206086
$ws.Cells.Item(7, 9).Value = 0.3235223881606854
$ws.Cells.Item(7, 10).Value = 0.3235223881606855
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.564779
$ws.Cells.Item(7, 14).Value = 1.694337
$ws.Cells.Item(7, 15).Value = 0.5632934075199009
$ws.Cells.Item(7, 16).Value = 0.5632934075199009
$ws.Cells.Item(7, 17).Value = 59.34059268166467
$ws.Cells.Item(7, 18).Value = 534.0653341349821
$ws.Cells.Item(7, 19).Value = 0.1822380284360085
$ws.Cells.Item(7, 20).Value = 0.1822380284360086
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Hspg2"
$ws.Cells.Item(8, 3).Value = "Col13a1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.2092613333333333
$ws.Cells.Item(8, 8).Value = 0.6277839999999999
$ws.Cells.Item(8, 9).Value = 0.0006443472634251982
$ws.Cells.Item(8, 10).Value = 0.0006443472634251983
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.3943663333333334
$ws.Cells.Item(8, 14).Value = 1.183099
$ws.Cells.Item(8, 15).Value = 0.393328993667368
$ws.Cells.Item(8, 16).Value = 0.393328993667368
$ws.Cells.Item(8, 17).Value = 0.0825256247351111
$ws.Cells.Item(8, 18).Value = 0.742730622616
$ws.Cells.Item(8, 19).Value = 0.0002534404606953556
$ws.Cells.Item(8, 20).Value = 0.0002534404606953557
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Hspg2"
$ws.Cells.Item(9, 3).Value = "Col13a1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.2092613333333333
$ws.Cells.Item(9, 8).Value = 0.6277839999999999
$ws.Cells.Item(9, 9).Value = 0.0006443472634251982
$ws.Cells.Item(9, 10).Value = 0.0006443472634251983
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.043492
$ws.Cells.Item(9, 14).Value = 0.130476
$ws.Cells.Item(9, 15).Value = 0.04337759881273123
$ws.Cells.Item(9, 16).Value = 0.04337759881273123
$ws.Cells.Item(9, 17).Value = 0.009101193909333333
$ws.Cells.Item(9, 18).Value = 0.08191074518399999
$ws.Cells.Item(9, 19).Value = 0.00002795023708893949
$ws.Cells.Item(9, 20).Value = 0.0000279502370889395
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Hspg2"
$ws.Cells.Item(10, 3).Value = "Col13a1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.2092613333333333
$ws.Cells.Item(10, 8).Value = 0.6277839999999999
$ws.Cells.Item(10, 9).Value = 0.0006443472634251982
$ws.Cells.Item(10, 10).Value = 0.0006443472634251983
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.564779
$ws.Cells.Item(10, 14).Value = 1.694337
$ws.Cells.Item(10, 15).Value = 0.5632934075199009
$ws.Cells.Item(10, 16).Value = 0.5632934075199009
$ws.Cells.Item(10, 17).Value = 0.1181864065786667
$ws.Cells.Item(10, 18).Value = 1.063677659208
$ws.Cells.Item(10, 19).Value = 0.0003629565656409031
$ws.Cells.Item(10, 20).Value = 0.0003629565656409031
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Hspg2"
$ws.Cells.Item(11, 3).Value = "Col13a1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.2582773333333333
$ws.Cells.Item(11, 8).Value = 0.774832
$ws.Cells.Item(11, 9).Value = 0.000795274933439325
$ws.Cells.Item(11, 10).Value = 0.0007952749334393251
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.3943663333333334
$ws.Cells.Item(11, 14).Value = 1.183099
$ws.Cells.Item(11, 15).Value = 0.393328993667368
$ws.Cells.Item(11, 16).Value = 0.393328993667368
$ws.Cells.Item(11, 17).Value = 0.1018558849297778
$ws.Cells.Item(11, 18).Value = 0.916702964368
$ws.Cells.Item(11, 19).Value = 0.0003128046892585727
$ws.Cells.Item(11, 20).Value = 0.0003128046892585728
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Hspg2"
$ws.Cells.Item(12, 3).Value = "Col13a1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.2582773333333333
$ws.Cells.Item(12, 8).Value = 0.774832
$ws.Cells.Item(12, 9).Value = 0.000795274933439325
$ws.Cells.Item(12, 10).Value = 0.0007952749334393251
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.043492
$ws.Cells.Item(12, 14).Value = 0.130476
$ws.Cells.Item(12, 15).Value = 0.04337759881273123
$ws.Cells.Item(12, 16).Value = 0.04337759881273123
$ws.Cells.Item(12, 17).Value = 0.01123299778133333
$ws.Cells.Item(12, 18).Value = 0.101096980032
$ws.Cells.Item(12, 19).Value = 0.00003449711700855256
$ws.Cells.Item(12, 20).Value = 0.00003449711700855257
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Hspg2"
$ws.Cells.Item(13, 3).Value = "Col13a1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.2582773333333333
$ws.Cells.Item(13, 8).Value = 0.774832
$ws.Cells.Item(13, 9).Value = 0.000795274933439325
$ws.Cells.Item(13, 10).Value = 0.0007952749334393251
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.564779
$ws.Cells.Item(13, 14).Value = 1.694337
$ws.Cells.Item(13, 15).Value = 0.5632934075199009
$ws.Cells.Item(13, 16).Value = 0.5632934075199009
$ws.Cells.Item(13, 17).Value = 0.1458696140426667
$ws.Cells.Item(13, 18).Value = 1.312826526384
$ws.Cells.Item(13, 19).Value = 0.0004479731271721997
$ws.Cells.Item(13, 20).Value = 0.0004479731271721998
$ws.Cells.Item(14, 1).Value = "Neutro"
$ws.Cells.Item(14, 2).Value = "Hspg2"
$ws.Cells.Item(14, 3).Value = "Col13a1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.4764636666666667
$ws.Cells.Item(14, 8).Value = 1.429391
$ws.Cells.Item(14, 9).Value = 0.001467103620376766
$ws.Cells.Item(14, 10).Value = 0.001467103620376766
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.3943663333333334
$ws.Cells.Item(14, 14).Value = 1.183099
$ws.Cells.Item(14, 15).Value = 0.393328993667368
$ws.Cells.Item(14, 16).Value = 0.393328993667368
$ws.Cells.Item(14, 17).Value = 0.1879012291898889
$ws.Cells.Item(14, 18).Value = 1.691111062709
$ws.Cells.Item(14, 19).Value = 0.0005770543906085456
$ws.Cells.Item(14, 20).Value = 0.0005770543906085456
$ws.Cells.Item(15, 1).Value = "Neutro"
$ws.Cells.Item(15, 2).Value = "Hspg2"
$ws.Cells.Item(15, 3).Value = "Col13a1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.4764636666666667
$ws.Cells.Item(15, 8).Value = 1.429391
$ws.Cells.Item(15, 9).Value = 0.001467103620376766
$ws.Cells.Item(15, 10).Value = 0.001467103620376766
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.043492
$ws.Cells.Item(15, 14).Value = 0.130476
$ws.Cells.Item(15, 15).Value = 0.04337759881273123
$ws.Cells.Item(15, 16).Value = 0.04337759881273123
$ws.Cells.Item(15, 17).Value = 0.02072235779066667
$ws.Cells.Item(15, 18).Value = 0.186501220116
$ws.Cells.Item(15, 19).Value = 0.00006363943226140888
$ws.Cells.Item(15, 20).Value = 0.0000636394322614089
$ws.Cells.Item(16, 1).Value = "Neutro"
$ws.Cells.Item(16, 2).Value = "Hspg2"
$ws.Cells.Item(16, 3).Value = "Col13a1"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.4764636666666667
$ws.Cells.Item(16, 8).Value = 1.429391
$ws.Cells.Item(16, 9).Value = 0.001467103620376766
$ws.Cells.Item(16, 10).Value = 0.001467103620376766
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.564779
$ws.Cells.Item(16, 14).Value = 1.694337
$ws.Cells.Item(16, 15).Value = 0.5632934075199009
$ws.Cells.Item(16, 16).Value = 0.5632934075199009
$ws.Cells.Item(16, 17).Value = 0.2690966731963333
$ws.Cells.Item(16, 18).Value = 2.421870058767
$ws.Cells.Item(16, 19).Value = 0.0008264097975068116
$ws.Cells.Item(16, 20).Value = 0.0008264097975068117
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Hspg2"
$ws.Cells.Item(17, 3).Value = "Col13a1"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 79.72331333333334
$ws.Cells.Item(17, 8).Value = 239.16994
$ws.Cells.Item(17, 9).Value = 0.2454801274523863
$ws.Cells.Item(17, 10).Value = 0.2454801274523863
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.3943663333333334
$ws.Cells.Item(17, 14).Value = 1.183099
$ws.Cells.Item(17, 15).Value = 0.393328993667368
$ws.Cells.Item(17, 16).Value = 0.393328993667368
$ws.Cells.Item(17, 17).Value = 31.44019076045112
$ws.Cells.Item(17, 18).Value = 282.96171684406
$ws.Cells.Item(17, 19).Value = 0.09655445149618433
$ws.Cells.Item(17, 20).Value = 0.09655445149618433
$ws.Cells.Item(18, 1).Value = "sCs"
$ws.Cells.Item(18, 2).Value = "Hspg2"
$ws.Cells.Item(18, 3).Value = "Col13a1"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 79.72331333333334
$ws.Cells.Item(18, 8).Value = 239.16994
$ws.Cells.Item(18, 9).Value = 0.2454801274523863
$ws.Cells.Item(18, 10).Value = 0.2454801274523863
$ws.Cells.Item(18, 11).Value = 1
$ws.Cells.Item(18, 12).Value = 0.3333333333333333
$ws.Cells.Item(18, 13).Value = 0.043492
$ws.Cells.Item(18, 14).Value = 0.130476
$ws.Cells.Item(18, 15).Value = 0.04337759881273123
$ws.Cells.Item(18, 16).Value = 0.04337759881273123
$ws.Cells.Item(18, 17).Value = 3.467326343493334
$ws.Cells.Item(18, 18).Value = 31.20593709144
$ws.Cells.Item(18, 19).Value = 0.01064833848512774
$ws.Cells.Item(18, 20).Value = 0.01064833848512774
$ws.Cells.Item(19, 1).Value = "sCs"
$ws.Cells.Item(19, 2).Value = "Hspg2"
$ws.Cells.Item(19, 3).Value = "Col13a1"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 79.72331333333334
$ws.Cells.Item(19, 8).Value = 239.16994
$ws.Cells.Item(19, 9).Value = 0.2454801274523863
$ws.Cells.Item(19, 10).Value = 0.2454801274523863
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.564779
$ws.Cells.Item(19, 14).Value = 1.694337
$ws.Cells.Item(19, 15).Value = 0.5632934075199009
$ws.Cells.Item(19, 16).Value = 0.5632934075199009
$ws.Cells.Item(19, 17).Value = 45.02605318108667
$ws.Cells.Item(19, 18).Value = 405.23447862978
$ws.Cells.Item(19, 19).Value = 0.1382773374710743
$ws.Cells.Item(19, 20).Value = 0.1382773374710743
